$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 with the same style as the other header cells (copy from E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Fill in the time_taken values for rows 2-24
$timestamps = @(
    "2021-10-05 10:50:17.185835",
    "2021-10-05 10:50:17.185846",
    "2021-10-05 10:50:17.185849",
    "2021-10-05 10:50:17.185851",
    "2021-10-05 10:50:17.185854",
    "2021-10-05 10:50:17.185857",
    "2021-10-05 10:50:17.185859",
    "2021-10-05 10:50:17.185862",
    "2021-10-05 10:50:17.185865",
    "2021-10-05 10:50:17.185868",
    "2021-10-05 10:50:17.185870",
    "2021-10-05 10:50:17.185872",
    "2021-10-05 10:50:17.185875",
    "2021-10-05 10:50:17.185877",
    "2021-10-05 10:50:17.185880",
    "2021-10-05 10:50:17.185882",
    "2021-10-05 10:50:17.185885",
    "2021-10-05 10:50:17.185888",
    "2021-10-05 10:50:17.185890",
    "2021-10-05 10:50:17.185893",
    "2021-10-05 10:50:17.185895",
    "2021-10-05 10:50:17.185898",
    "2021-10-05 10:50:17.185900"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

Write-Output "done"
